$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText-NoBoldChange($cellRange, $oldText, $newText) {
    $find = $cellRange.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $oldText
    $find.Replacement.Text = $newText
    $find.Forward = $true
    $find.Wrap = 0
    $find.Format = $false
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.Execute([ref]$null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null
}

# 1. Header cell: "Characteristic" -> "Baseline Characteristics", un-bold
$headerCellRange = $t.Rows.Item(1).Cells.Item(1).Range
$find = $headerCellRange.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Characteristic"
$find.Replacement.Text = "Baseline Characteristics"
$find.Replacement.Font.Bold = 0
$find.Forward = $true
$find.Wrap = 0
$find.Format = $true
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.Execute([ref]$null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null

# 2. "Nonwhite" -> "Non-white"
$d.Content.Find.Execute("Nonwhite", $true, $false, $false, $false, $false, $true, 1, $false, "Non-white", 2) | Out-Null

# 3. Row 73 (1-indexed): "Greenspace 300m from residence, Median (Q1, Q3)" -> "water_300m, Median (Q1, Q3)"
#    and row height 624 -> 630 (twips == 31.2pt -> 31.5pt)
$row73 = $t.Rows.Item(73)
Set-CellText-NoBoldChange $row73.Cells.Item(1).Range "Greenspace 300m from residence, Median (Q1, Q3)" "water_300m, Median (Q1, Q3)"
$row73.Height = 31.5
$row73.HeightRule = 0

# 4. Row 74 (1-indexed): "Water 300m from residence, Median (Q1, Q3)" -> "greenspace_300m, Median (Q1, Q3)"
#    and row height 607 -> 630 (twips == 30.35pt -> 31.5pt)
$row74 = $t.Rows.Item(74)
Set-CellText-NoBoldChange $row74.Cells.Item(1).Range "Water 300m from residence, Median (Q1, Q3)" "greenspace_300m, Median (Q1, Q3)"
$row74.Height = 31.5
$row74.HeightRule = 0

# 5. Row 75 (1-indexed): "Natural env. 300m from residence, Median (Q1, Q3)" -> "naturalenv_300m, Median (Q1, Q3)"
#    and row height 607 -> 630 (twips == 30.35pt -> 31.5pt)
$row75 = $t.Rows.Item(75)
Set-CellText-NoBoldChange $row75.Cells.Item(1).Range "Natural env. 300m from residence, Median (Q1, Q3)" "naturalenv_300m, Median (Q1, Q3)"
$row75.Height = 31.5
$row75.HeightRule = 0
